# Resolve attendee names to IDs for accurate seating connections.
#
# The "Contestants" rows got rewritten (their A:M data now appears one row
# "earlier" than before, wrapping around), the "Seat Assignments" sheet
# picked up a new seat-B1 row for Kathleen Reynolds (and the existing row's
# seat moved from A1 to B2), and the "Standbys" sheet lost the row that
# used to hold Kathleen's pending-standby entry, leaving only Peter's.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Contestants sheet — re-map each row's data.
# ---------------------------------------------------------------------
$contestants = $wb.Worksheets.Item("Contestants")

# Row 2 -> Peter Adamidis
$contestants.Range("A2").Value = "0ccaf8bc-6ade-4ad6-9537-92f9829a6502"
$contestants.Range("B2").Value = "Peter Adamidis"
$contestants.Range("C2").Value = 34
$contestants.Range("E2").Value = "peter.adamidis@gmail.com"
$contestants.Range("G2").Value = ""
$contestants.Range("H2").ClearContents()
$contestants.Range("J2").Value = "Kathleen Reynolds, Felicity Parker-Hill"
$contestants.Range("L2").Value = "Y"
$contestants.Range("M2").Value = "Broken Leg"

# Row 3 -> Kathleen Reynolds
$contestants.Range("A3").Value = "d698b1de-6641-45c6-aa63-f577d2b634bb"
$contestants.Range("B3").Value = "Kathleen Reynolds"
$contestants.Range("C3").Value = 33
$contestants.Range("E3").Value = "kathleenmonicareynolds@gmail.com"
$contestants.Range("G3").Value = "Footscray"
$contestants.Range("H3").ClearContents()
$contestants.Range("J3").Value = "Peter Adamidis, Felicity Parker-Hill"
$contestants.Range("L3").Value = "N"
$contestants.Range("M3").Value = "N/A"

# Row 4 -> Felicity Parker-Hill
$contestants.Range("A4").Value = "28603f95-d5f6-47ab-88c4-0d79742a6b02"
$contestants.Range("B4").Value = "Felicity Parker-Hill"
$contestants.Range("C4").Value = 27
$contestants.Range("E4").Value = "felicity.parkerhill@endemolshine.com.au"
$contestants.Range("G4").Value = "Melbourne"
$contestants.Range("H4").Value = ""
$contestants.Range("J4").Value = "Peter Adamidis, Kathleen Reynolds"
$contestants.Range("L4").Value = "N"
$contestants.Range("M4").Value = "N/A"

# ---------------------------------------------------------------------
# 2) Seat Assignments sheet — update seat for the existing row and add a
#    newly-resolved seat assignment for Kathleen Reynolds.
# ---------------------------------------------------------------------
$seats = $wb.Worksheets.Item("Seat Assignments")

$seats.Range("A2").Value = "91e25164-6f67-42f7-b978-9132a406c060"
$seats.Range("E2").Value = "B2"

$seats.Range("A3").Value = "ff87f03b-8891-4bb6-ac5c-a510d216fdd6"
$seats.Range("B3").Value = "e432f0fe-1383-44a2-990c-5f787da5008a"
$seats.Range("C3").Value = "d698b1de-6641-45c6-aa63-f577d2b634bb"
$seats.Range("D3").Value = 1
$seats.Range("E3").Value = "B1"

# ---------------------------------------------------------------------
# 3) Standbys sheet — Kathleen Reynolds' standby row is resolved away,
#    leaving only Peter Adamidis' pending standby (now on row 2).
# ---------------------------------------------------------------------
$standbys = $wb.Worksheets.Item("Standbys")

$standbys.Range("A2").Value = "69511596-e3b5-41be-a93b-920748af4fe0"
$standbys.Range("C2").Value = "0ccaf8bc-6ade-4ad6-9537-92f9829a6502"

$standbys.Rows(3).Delete()
